# Scheduled market-data refresh: update recomputed price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the leve-profit
# sheets, per the latest Sargatanas market board scrape.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1900.2858
$ws.Range("I11").Value = 1900.2858
$ws.Range("K11").Value = 1900.2858
$ws.Range("M11").Value = -1760.2858

$ws.Range("H106").Value = 38463344
$ws.Range("J106").Value = 1598.5
$ws.Range("L106").Value = 1598.5
$ws.Range("N106").Value = -2860.5

$ws.Range("H138").Value = 5272197.5
$ws.Range("I138").Value = 2697.4285
$ws.Range("J138").Value = 8346072.5
$ws.Range("K138").Value = 8092.2855
$ws.Range("L138").Value = 25038217.5
$ws.Range("M138").Value = -2952.2855
$ws.Range("N138").Value = -25048497.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 58825508
$ws.Range("I2").Value = 1452.1666
$ws.Range("J2").Value = 200003250
$ws.Range("K2").Value = 1452.1666
$ws.Range("L2").Value = 200003250
$ws.Range("M2").Value = -1339.1666
$ws.Range("N2").Value = -200003476

$ws.Range("H6").Value = 20000002
$ws.Range("I6").Value = 20000002
$ws.Range("K6").Value = 20000002
$ws.Range("M6").Value = -19999829

$ws.Range("H102").Value = 10528576
$ws.Range("I102").Value = 14287710
$ws.Range("K102").Value = 14287710
$ws.Range("M102").Value = -14286088

$ws.Range("H116").Value = 58825508
$ws.Range("I116").Value = 1452.1666
$ws.Range("J116").Value = 200003250
$ws.Range("K116").Value = 1452.1666
$ws.Range("L116").Value = 200003250
$ws.Range("M116").Value = 841.8334
$ws.Range("N116").Value = -200007838

$ws.Range("H132").Value = 4028.3035
$ws.Range("I132").Value = 2231.6462
$ws.Range("K132").Value = 6694.9386
$ws.Range("M132").Value = -4164.9386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 58825508
$ws.Range("I3").Value = 1452.1666
$ws.Range("J3").Value = 200003250
$ws.Range("K3").Value = 1452.1666
$ws.Range("L3").Value = 200003250
$ws.Range("M3").Value = -1338.1666
$ws.Range("N3").Value = -200003478

$ws.Range("H20").Value = 10420238
$ws.Range("I20").Value = 18521192
$ws.Range("J20").Value = 4726.7144
$ws.Range("K20").Value = 18521192
$ws.Range("L20").Value = 4726.7144
$ws.Range("M20").Value = -18520945
$ws.Range("N20").Value = -5220.7144

$ws.Range("H80").Value = 31250492
$ws.Range("I80").Value = 55556056
$ws.Range("J80").Value = 482.7143
$ws.Range("K80").Value = 55556056
$ws.Range("L80").Value = 482.7143
$ws.Range("M80").Value = -55555058
$ws.Range("N80").Value = -2478.7143

$ws.Range("H83").Value = 31250492
$ws.Range("I83").Value = 55556056
$ws.Range("J83").Value = 482.7143
$ws.Range("K83").Value = 277780280
$ws.Range("L83").Value = 2413.5715
$ws.Range("M83").Value = -277775288
$ws.Range("N83").Value = -12397.5715

$ws.Range("H94").Value = 2137.0527
$ws.Range("I94").Value = 1390.6428
$ws.Range("K94").Value = 1390.6428
$ws.Range("M94").Value = -939.6428000000001

$ws.Range("H105").Value = 2558
$ws.Range("I105").Value = 1492.8948
$ws.Range("K105").Value = 1492.8948
$ws.Range("M105").Value = 254.1052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 352.5
$ws.Range("I11").Value = 505
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 505
$ws.Range("L11").Value = 200
$ws.Range("M11").Value = -365
$ws.Range("N11").Value = -480

$ws.Range("H19").Value = 2525.6
$ws.Range("I19").Value = 2964.75
$ws.Range("K19").Value = 2964.75
$ws.Range("M19").Value = -2794.75

$ws.Range("H24").Value = 2525.6
$ws.Range("I24").Value = 2964.75
$ws.Range("K24").Value = 2964.75
$ws.Range("M24").Value = -2794.75

$ws.Range("H31").Value = 7015.911
$ws.Range("I31").Value = 2766.5
$ws.Range("K31").Value = 2766.5
$ws.Range("M31").Value = -2471.5

$ws.Range("H34").Value = 7015.911
$ws.Range("I34").Value = 2766.5
$ws.Range("K34").Value = 2766.5
$ws.Range("M34").Value = -2564.5

$ws.Range("H99").Value = 12001.625
$ws.Range("I99").Value = 14599.8
$ws.Range("K99").Value = 14599.8
$ws.Range("M99").Value = -13101.8

$ws.Range("H126").Value = 12001.625
$ws.Range("I126").Value = 14599.8
$ws.Range("K126").Value = 43799.39999999999
$ws.Range("M126").Value = -41329.39999999999

$ws.Range("H132").Value = 6776.143
$ws.Range("I132").Value = 3836.1667
$ws.Range("K132").Value = 11508.5001
$ws.Range("M132").Value = -8978.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 5890523
$ws.Range("I55").Value = 1989
$ws.Range("J55").Value = 7152352
$ws.Range("K55").Value = 5967
$ws.Range("L55").Value = 21457056
$ws.Range("M55").Value = -5790
$ws.Range("N55").Value = -21457410

$ws.Range("H98").Value = 2583.5454
$ws.Range("J98").Value = 3970.1428
$ws.Range("L98").Value = 11910.4284
$ws.Range("N98").Value = -14906.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6057139.5
$ws.Range("I122").Value = 18158392
$ws.Range("K122").Value = 54475176
$ws.Range("M122").Value = -54472726

$ws.Range("H132").Value = 4481.0835
$ws.Range("I132").Value = 2546.5
$ws.Range("K132").Value = 7639.5
$ws.Range("M132").Value = -5109.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4556.4375
$ws.Range("I7").Value = 2818.9092
$ws.Range("K7").Value = 2818.9092
$ws.Range("M7").Value = -2706.9092

$ws.Range("H19").Value = 1001.5
$ws.Range("I19").Value = 1003
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 1003
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = -833
$ws.Range("N19").Value = -1340

$ws.Range("H55").Value = 360.25
$ws.Range("I55").Value = 24.88889
$ws.Range("J55").Value = 561.4666999999999
$ws.Range("K55").Value = 24.88889
$ws.Range("L55").Value = 561.4666999999999
$ws.Range("M55").Value = 148.11111
$ws.Range("N55").Value = -907.4666999999999

$ws.Range("H126").Value = 4556.4375
$ws.Range("I126").Value = 2818.9092
$ws.Range("K126").Value = 8456.7276
$ws.Range("M126").Value = -5986.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").Value = $null

$ws.Range("H126").Value = 3325.625
$ws.Range("I126").Value = 5900
$ws.Range("J126").Value = 2467.5
$ws.Range("K126").Value = 17700
$ws.Range("L126").Value = 7402.5
$ws.Range("M126").Value = -15230
$ws.Range("N126").Value = -12342.5

$ws.Range("H132").Value = 32264400
$ws.Range("I132").Value = 47626060
$ws.Range("K132").Value = 142878180
$ws.Range("M132").Value = -142875650
